# Update FY financial figures in the WMMVY yearly financials sheet
# (author commit message: "Doing Updates for Financials")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: Total Revenue
$ws.Range("D8").Value = 29649300
$ws.Range("E8").Value = 27534900
$ws.Range("F8").Value = 24614100
$ws.Range("G8").Value = 22807900
$ws.Range("H8").Value = 21920100
$ws.Range("I8").Value = 21311800
$ws.Range("J8").Value = 19594100

# Row 9: Cost of Revenue
$ws.Range("D9").Value = 23044900
$ws.Range("E9").Value = 21458600
$ws.Range("F9").Value = 19307500
$ws.Range("G9").Value = 17759000
$ws.Range("H9").Value = 17112800
$ws.Range("I9").Value = 16645200
$ws.Range("J9").Value = 15192200

# Row 10: Gross Profit
$ws.Range("D10").Value = 6604400
$ws.Range("E10").Value = 6076300
$ws.Range("F10").Value = 5306600
$ws.Range("G10").Value = 5048900
$ws.Range("H10").Value = 4807300
$ws.Range("I10").Value = 4666600
$ws.Range("J10").Value = 4401900

# Row 17: Total Operating Expenses
$ws.Range("D17").Value = 27382000
$ws.Range("E17").Value = 25494300
$ws.Range("F17").Value = 22916200
$ws.Range("G17").Value = 21012400
$ws.Range("H17").Value = 20283900
$ws.Range("I17").Value = 19686600
$ws.Range("J17").Value = 18066600

# Row 18: Operating Income or Loss
$ws.Range("D18").Value = 2267300
$ws.Range("E18").Value = 2040600
$ws.Range("F18").Value = 1697900
$ws.Range("G18").Value = 1795500
$ws.Range("H18").Value = 1636200
$ws.Range("I18").Value = 1625100
$ws.Range("J18").Value = 1527500

# Row 20: Total Other Income/Expenses Net
$ws.Range("D20").Value = 46000
$ws.Range("E20").Value = 49000
$ws.Range("F20").Value = 127200
$ws.Range("G20").Value = 54700
$ws.Range("H20").Value = 56800
$ws.Range("I20").Value = 128400
$ws.Range("J20").Value = 59300

# Row 21: Earnings Before Interest And Taxes
$ws.Range("D21").Value = 2915000
$ws.Range("E21").Value = 2642200
$ws.Range("F21").Value = 2343100
$ws.Range("G21").Value = 2329700
$ws.Range("H21").Value = 2142000
$ws.Range("I21").Value = 2187200
$ws.Range("J21").Value = "NA"

# Row 22: Interest Expense
$ws.Range("D22").Value = 74300
$ws.Range("E22").Value = 65700
$ws.Range("F22").Value = 122600
$ws.Range("G22").Value = 62700
$ws.Range("H22").Value = 57600
$ws.Range("I22").Value = 107700
$ws.Range("J22").Value = 46600

# Row 23: Income Before Tax
$ws.Range("D23").Value = 2239000
$ws.Range("E23").Value = 2023900
$ws.Range("F23").Value = 1702500
$ws.Range("G23").Value = 1787600
$ws.Range("H23").Value = 1635400
$ws.Range("I23").Value = 1645900
$ws.Range("J23").Value = 1540200

# Row 24: Income Tax Expense
$ws.Range("D24").Value = 563700
$ws.Range("E24").Value = 549400
$ws.Range("F24").Value = 489900
$ws.Range("G24").Value = 492400
$ws.Range("H24").Value = 498100
$ws.Range("I24").Value = 477700
$ws.Range("J24").Value = 398000

# Row 26: Income After Tax
$ws.Range("D26").Value = 1675200
$ws.Range("E26").Value = 1474500
$ws.Range("F26").Value = 1212500
$ws.Range("G26").Value = 1295100
$ws.Range("H26").Value = 1137300
$ws.Range("I26").Value = 1168100
$ws.Range("J26").Value = 1142200

# Row 27: Net Income From Continuing Ops
$ws.Range("D27").Value = 1675200
$ws.Range("E27").Value = 1474500
$ws.Range("F27").Value = 1212400
$ws.Range("G27").Value = 1294600
$ws.Range("H27").Value = 1137500
$ws.Range("I27").Value = 1168500
$ws.Range("J27").Value = 1142000

# Row 29: Discontinued Operations
$ws.Range("D29").Value = 386600
$ws.Range("E29").Value = 250500
$ws.Range("F29").Value = 151800
$ws.Range("G29").Value = 279000
$ws.Range("H29").Value = 37500
$ws.Range("I29").Value = 35300
$ws.Range("J29").Value = "NA"

# Row 32: Other Items
$ws.Range("D32").Value = -46000
$ws.Range("E32").Value = -49000
$ws.Range("F32").Value = -127200
$ws.Range("G32").Value = -54700
$ws.Range("H32").Value = -56800
$ws.Range("I32").Value = -128400
$ws.Range("J32").Value = -59300

# Row 33: Net Income
$ws.Range("D33").Value = 2061800
$ws.Range("E33").Value = 1725000
$ws.Range("F33").Value = 1364200
$ws.Range("G33").Value = 1573600
$ws.Range("H33").Value = 1174900
$ws.Range("I33").Value = 1203800
$ws.Range("J33").Value = 1142000

# Row 35: Net Income Applicable To Common Shares
$ws.Range("D35").Value = 2061800
$ws.Range("E35").Value = 1725000
$ws.Range("F35").Value = 1364200
$ws.Range("G35").Value = 1573600
$ws.Range("H35").Value = 1174900
$ws.Range("I35").Value = 1203800
$ws.Range("J35").Value = 1142000

# Row 41: Cash And Cash Equivalents
$ws.Range("D41").Value = 1841000
$ws.Range("E41").Value = 1446900
$ws.Range("F41").Value = 1282200
$ws.Range("G41").Value = 2901300
$ws.Range("H41").Value = 2185600
$ws.Range("I41").Value = 1456600
$ws.Range("J41").Value = 2603200

# Row 43: Net Receivables
$ws.Range("D43").Value = 560400
$ws.Range("E43").Value = 507800
$ws.Range("F43").Value = 504700
$ws.Range("G43").Value = 570400
$ws.Range("H43").Value = 656900
$ws.Range("I43").Value = 536700
$ws.Range("J43").Value = 593800

# Row 44: Inventory
$ws.Range("D44").Value = 3075400
$ws.Range("E44").Value = 2775600
$ws.Range("F44").Value = 2573000
$ws.Range("G44").Value = 4879800
$ws.Range("H44").Value = 4530100
$ws.Range("I44").Value = 2021800
$ws.Range("J44").Value = 2141200

# Row 45: Other Current Assets
$ws.Range("D45").Value = 50800
$ws.Range("E45").Value = 649800
$ws.Range("F45").Value = 28500
$ws.Range("G45").Value = 464200
$ws.Range("H45").Value = 267100
$ws.Range("I45").Value = 52600
$ws.Range("J45").Value = 48800

# Row 46: Total Current Assets
$ws.Range("D46").Value = 5527600
$ws.Range("E46").Value = 5380100
$ws.Range("F46").Value = 4388400
$ws.Range("G46").Value = 4678200
$ws.Range("H46").Value = 4274600
$ws.Range("I46").Value = 4067600
$ws.Range("J46").Value = 3948700

# Row 47: Long Term Investments
$ws.Range("D47").Value = "NA"
$ws.Range("E47").Value = "NA"
$ws.Range("F47").Value = "NA"
$ws.Range("G47").Value = 3800
$ws.Range("H47").Value = 10700
$ws.Range("I47").Value = "NA"
$ws.Range("J47").Value = 15300

# Row 48: Property Plant and Equipment
$ws.Range("D48").Value = 7245000
$ws.Range("E48").Value = 7052000
$ws.Range("F48").Value = 6735100
$ws.Range("G48").Value = 6047900
$ws.Range("H48").Value = 5664400
$ws.Range("I48").Value = 6070700
$ws.Range("J48").Value = 8025400

# Row 49: Goodwill
$ws.Range("D49").Value = 2036100
$ws.Range("E49").Value = 2138100
$ws.Range("F49").Value = 1782100
$ws.Range("G49").Value = 3014100
$ws.Range("H49").Value = 1342500
$ws.Range("I49").Value = 1341000
$ws.Range("J49").Value = 1645500

# Row 52: Other Assets
$ws.Range("D52").Value = 461900
$ws.Range("E52").Value = 340300
$ws.Range("F52").Value = 213200
$ws.Range("G52").Value = 185900
$ws.Range("H52").Value = 48600
$ws.Range("I52").Value = 17200
$ws.Range("J52").Value = 26300

# Row 54: Total Assets
$ws.Range("D54").Value = 15270600
$ws.Range("E54").Value = 14910400
$ws.Range("F54").Value = 13118800
$ws.Range("G54").Value = 12863500
$ws.Range("H54").Value = 11909200
$ws.Range("I54").Value = 11496600
$ws.Range("J54").Value = 11338100

# Row 57: Accounts Payable
$ws.Range("D57").Value = 4142700
$ws.Range("E57").Value = 3409300
$ws.Range("F57").Value = 2916800
$ws.Range("G57").Value = 2726200
$ws.Range("H57").Value = 2462400
$ws.Range("I57").Value = 2315500
$ws.Range("J57").Value = 2630200

# Row 58: Short/Current Long Term Debt
$ws.Range("D58").Value = 19700
$ws.Range("E58").Value = 27600
$ws.Range("F58").Value = 28800
$ws.Range("G58").Value = 28900
$ws.Range("H58").Value = 37400
$ws.Range("I58").Value = 24600
$ws.Range("J58").Value = 20800

# Row 59: Other Current Liabilities
$ws.Range("D59").Value = 1414400
$ws.Range("E59").Value = 2705700
$ws.Range("F59").Value = 1052900
$ws.Range("G59").Value = 2405500
$ws.Range("H59").Value = 1710800
$ws.Range("I59").Value = 875000
$ws.Range("J59").Value = 914900

# Row 60: Total Current Liabilities
$ws.Range("D60").Value = 5576800
$ws.Range("E60").Value = 4878100
$ws.Range("F60").Value = 3998400
$ws.Range("G60").Value = 3906300
$ws.Range("H60").Value = 3399100
$ws.Range("I60").Value = 3215100
$ws.Range("J60").Value = 3443300

# Row 61: Long Term Debt
$ws.Range("D61").Value = 709900
$ws.Range("E61").Value = 676300
$ws.Range("F61").Value = 677700
$ws.Range("G61").Value = 681700
$ws.Range("H61").Value = 712000
$ws.Range("I61").Value = 653700
$ws.Range("J61").Value = 848600

# Row 62: Other Liabilities
$ws.Range("D62").Value = 732800
$ws.Range("E62").Value = 717400
$ws.Range("F62").Value = 590400
$ws.Range("G62").Value = 536400
$ws.Range("H62").Value = 431500
$ws.Range("I62").Value = 395500
$ws.Range("J62").Value = 1212400

# Row 66: Total Liabilities
$ws.Range("D66").Value = 7019500
$ws.Range("E66").Value = 6271900
$ws.Range("F66").Value = 5268000
$ws.Range("G66").Value = 5094000
$ws.Range("H66").Value = 4516800
$ws.Range("I66").Value = 4271300
$ws.Range("J66").Value = 4673100

# Row 72: Retained Earnings
$ws.Range("D72").Value = 4680400
$ws.Range("E72").Value = 5336500
$ws.Range("F72").Value = 4928600
$ws.Range("G72").Value = 10296500
$ws.Range("H72").Value = 6440800
$ws.Range("I72").Value = 5088600
$ws.Range("J72").Value = 4504300

# Row 76: Total Stockholder Equity
$ws.Range("D76").Value = 8251100
$ws.Range("E76").Value = 8638600
$ws.Range("F76").Value = 7850800
$ws.Range("G76").Value = 7769500
$ws.Range("H76").Value = 7392400
$ws.Range("I76").Value = 7225300
$ws.Range("J76").Value = 6665000

# Row 81: Net Income
$ws.Range("D81").Value = 2061800
$ws.Range("E81").Value = 1725000
$ws.Range("F81").Value = 1364200
$ws.Range("G81").Value = 1573600
$ws.Range("H81").Value = 1174900
$ws.Range("I81").Value = 1203800
$ws.Range("J81").Value = 1142000

# Row 83: Depreciation
$ws.Range("D83").Value = 602300
$ws.Range("E83").Value = 553100
$ws.Range("F83").Value = 518500
$ws.Range("G83").Value = 479900
$ws.Range("H83").Value = 449400
$ws.Range("I83").Value = 434100
$ws.Range("J83").Value = "NA"

# Row 89: Total Cash Flow From Operating Activities
$ws.Range("D89").Value = 2541300
$ws.Range("E89").Value = 2436400
$ws.Range("F89").Value = 2012100
$ws.Range("G89").Value = 1973200
$ws.Range("H89").Value = 1484500
$ws.Range("I89").Value = 1533000
$ws.Range("J89").Value = 1664800

# Row 91: Capital Expenditures
$ws.Range("D91").Value = -901300
$ws.Range("E91").Value = -741400
$ws.Range("F91").Value = -647900
$ws.Range("G91").Value = -656400
$ws.Range("H91").Value = -723400
$ws.Range("I91").Value = -758200
$ws.Range("J91").Value = -949200

# Row 94: Total Cash Flows From Investing Activities
$ws.Range("D94").Value = 201000
$ws.Range("E94").Value = -713800
$ws.Range("F94").Value = -415200
$ws.Range("G94").Value = -212700
$ws.Range("H94").Value = -739100
$ws.Range("I94").Value = -739800
$ws.Range("J94").Value = "NA"

# Row 96: Dividends Paid
$ws.Range("D96").Value = -2211300
$ws.Range("E96").Value = -1498400
$ws.Range("F96").Value = -1632400
$ws.Range("G96").Value = -1119400
$ws.Range("H96").Value = -830400
$ws.Range("I96").Value = -497100
$ws.Range("J96").Value = -499600

# Row 100: Total Cash Flows From Financing Activities
$ws.Range("D100").Value = -2295400
$ws.Range("E100").Value = -1582300
$ws.Range("F100").Value = -1810800
$ws.Range("G100").Value = -1409000
$ws.Range("H100").Value = -1074600
$ws.Range("I100").Value = -620200
$ws.Range("J100").Value = "NA"

# Row 101: Effect Of Exchange Rate Changes 
$ws.Range("D101").Value = -52800
$ws.Range("E101").Value = 24300
$ws.Range("F101").Value = 45400
$ws.Range("G101").Value = 6300
$ws.Range("H101").Value = -34600
$ws.Range("I101").Value = -18000
$ws.Range("J101").Value = "NA"

# Row 102: Change In Cash and Cash Equivalents 
$ws.Range("D102").Value = 394100
$ws.Range("E102").Value = 164700
$ws.Range("F102").Value = -168500
$ws.Range("G102").Value = 357800
$ws.Range("H102").Value = -363800
$ws.Range("I102").Value = 155000
$ws.Range("J102").Value = 26100
